# Update main GSC export data on the "Chart" sheet.
# The exported coverage data rolled forward by one day: the oldest day
# (2025-10-12) drops off, every remaining row shifts up to take the
# values previously belonging to the next day, and the series now ends
# one day earlier (2026-01-02) since no new trailing day was appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Deleting the first data row (row 2, date 2025-10-12) shifts every
# subsequent row up by one, which reproduces the new Date/Not
# indexed/Indexed/Impressions values for rows 2-84 and removes the
# former last row (2026-01-03) automatically.
$ws.Rows.Item(2).Delete()
